$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.926300000000003
$ws.Range("D2").Value = -7.350600000000001
$ws.Range("A3").Value = -21.9745
$ws.Range("C3").Value = -11.2788
$ws.Range("D6").Value = -7.670099999999996
$ws.Range("C12").Value = -11.81029999999999
$ws.Range("A14").Value = -21.6214
$ws.Range("A16").Value = -21.8525
$ws.Range("B18").Value = 6.551299999999994
$ws.Range("D19").Value = -9.043599999999991
$ws.Range("A21").Value = -20.51619999999999
$ws.Range("A23").Value = -21.01039999999998
$ws.Range("B24").Value = 6.921800000000003
$ws.Range("C24").Value = -12.60229999999999
$ws.Range("D24").Value = -7.531599999999997
$ws.Range("A25").Value = -21.54869999999998
$ws.Range("B25").Value = 5.8728
$ws.Range("C25").Value = -12.8883
$ws.Range("A26").Value = -21.13729999999997
$ws.Range("B27").Value = 5.734300000000002
$ws.Range("D27").Value = -8.654300000000001
$ws.Range("A29").Value = -20.99099999999997
$ws.Range("B30").Value = 5.792199999999997
$ws.Range("D30").Value = -7.160800000000003
$ws.Range("B31").Value = 5.158000000000002
$ws.Range("D31").Value = -8.537100000000002
$ws.Range("D33").Value = -7.913099999999998
$ws.Range("B39").Value = 9.879000000000001
$ws.Range("A40").Value = -20.35219999999999
$ws.Range("C41").Value = -13.2763
$ws.Range("B42").Value = 9.748499999999991
$ws.Range("D42").Value = -8.876299999999992
$ws.Range("B48").Value = 5.380900000000002
$ws.Range("C50").Value = -13.19979999999999
$ws.Range("B51").Value = 5.769500000000001
$ws.Range("B52").Value = 5.123300000000001
$ws.Range("A53").Value = -21.5186
$ws.Range("C53").Value = -10.4739
$ws.Range("B55").Value = 6.465299999999993
$ws.Range("D55").Value = -7.653499999999999
$ws.Range("B56").Value = 5.347500000000001
$ws.Range("C56").Value = -12.21239999999999
$ws.Range("A57").Value = -22.0862
$ws.Range("B57").Value = 5.476299999999996
$ws.Range("C57").Value = -12.94889999999999
$ws.Range("C58").Value = -13.79819999999999
$ws.Range("D58").Value = -8.177399999999999
$ws.Range("A59").Value = -22.2563
$ws.Range("B60").Value = 5.6502
$ws.Range("C61").Value = -12.66579999999999
$ws.Range("C63").Value = -11.4553
$ws.Range("C64").Value = -11.5006
$ws.Range("A65").Value = -21.91729999999999
$ws.Range("D65").Value = -8.0425
$ws.Range("A69").Value = -21.58929999999998
$ws.Range("C70").Value = -11.6935
$ws.Range("D70").Value = -8.251200000000001
$ws.Range("C72").Value = -11.3548
$ws.Range("B73").Value = 8.320099999999995
$ws.Range("B74").Value = 9.716099999999988
$ws.Range("D74").Value = -8.864899999999997
$ws.Range("D75").Value = -8.123099999999999
$ws.Range("A79").Value = -20.44360000000002
$ws.Range("A83").Value = -22.08329999999999
$ws.Range("D83").Value = -8.273399999999992
$ws.Range("D84").Value = -8.871499999999999
$ws.Range("C86").Value = -13.1047
$ws.Range("D86").Value = -8.395799999999991
$ws.Range("B89").Value = 5.761399999999996
$ws.Range("C89").Value = -10.4578
$ws.Range("B90").Value = 5.583700000000002
$ws.Range("B92").Value = 5.544999999999995
$ws.Range("A93").Value = -21.20739999999998
$ws.Range("D96").Value = -7.326499999999998
$ws.Range("D97").Value = -8.089499999999999
$ws.Range("C98").Value = -12.27309999999999
$ws.Range("A100").Value = -21.59359999999999
$ws.Range("C100").Value = -13.23699999999999
$ws.Range("C102").Value = -13.3136
